$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly data row above the current first data row (row 7),
# pushing the existing rows 7:19 down to 8:20.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with this week's observation
# (same market/category/quality as the rest of the series).
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(7, 3).Value = "Los Lagos"
$ws.Cells.Item(7, 4).Value = 44803
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 100112035
$ws.Cells.Item(7, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 90
$ws.Cells.Item(7, 11).Value = 24000
$ws.Cells.Item(7, 12).Value = 24000
$ws.Cells.Item(7, 13).Value = 24000
$ws.Cells.Item(7, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(7, 16).Value = 1600
$ws.Cells.Item(7, 17).Value = 15
$ws.Cells.Item(7, 18).Value = "Hortaliza"
